$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "43.323.71"
Set-TextCell "E2" "  -1.11%  "
Set-TextCell "D3" "2.350.76"
Set-TextCell "E3" "  +5.40%  "
Set-TextCell "E4" "  +0.03%  "
Set-TextCell "D5" "232.43"
Set-TextCell "E5" "  +0.93%  "
Set-TextCell "D6" "0.649"
Set-TextCell "E6" "  +0.78%  "
Set-TextCell "D7" "67.65"
Set-TextCell "E7" "  +6.87%  "
Set-TextCell "E8" "  +0.01%  "
Set-TextCell "E9" "  +1.81%  "
Set-TextCell "D11" "56.79"
Set-TextCell "D12" "26.29"
Set-TextCell "E12" "  -1.02%  "
Set-TextCell "D13" "2.702.00"
Set-TextCell "E13" "  +5.46%  "
Set-TextCell "E14" "  -0.78%  "
Set-TextCell "D15" "15.63"
Set-TextCell "E15" "  +1.56%  "
Set-TextCell "E16" "  +1.85%  "
Set-TextCell "E17" "  +1.76%  "
Set-TextCell "D18" "2.351.61"
Set-TextCell "E18" "  +5.83%  "
Set-TextCell "D19" "43.255.37"
Set-TextCell "E19" "  -0.76%  "
Set-TextCell "D20" "0.0₃0978"
Set-TextCell "E20" "  -0.43%  "
Set-TextCell "D21" "73.93"
Set-TextCell "E21" "  +1.82%  "
Set-TextCell "D22" "6.24"
Set-TextCell "E22" "  +3.72%  "
Set-TextCell "D23" "248.67"
Set-TextCell "E23" "  +0.11%  "
Set-TextCell "D24" "3.98"
Set-TextCell "E24" "  +17.35%  "
Set-TextCell "E25" "  -0.06%  "
Set-TextCell "E26" "  +1.23%  "
Set-TextCell "E27" "  -3.23%  "
Set-TextCell "D28" "9.89"
Set-TextCell "E28" "  +0.56%  "
Set-TextCell "D29" "22.31"
Set-TextCell "E29" "  +7.34%  "
Set-TextCell "D30" "172.62"
Set-TextCell "E30" "  +1.24%  "
Set-TextCell "D31" "1.54"
Set-TextCell "E31" "  +11.87%  "
Set-TextCell "E32" "  -7.27%  "
Set-TextCell "E33" "  +0.36%  "
Set-TextCell "D34" "5.02"
Set-TextCell "E34" "  +6.03%  "
Set-TextCell "D35" "0.0692"
Set-TextCell "E35" "  -0.57%  "
Set-TextCell "E36" "  +3.79%  "
Set-TextCell "D37" "2.49"
Set-TextCell "E37" "  +10.18%  "
Set-TextCell "E38" "  +2.08%  "
Set-TextCell "D39" "3.62"
Set-TextCell "E39" "  -0.51%  "
Set-TextCell "E40" "  -1.59%  "
Set-TextCell "B41" "FraxShare"
Set-TextCell "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D41" "8.95"
Set-TextCell "E41" "  +9.12%  "
Set-TextCell "B42" "BinanceUSD"
Set-TextCell "C42" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D42" "1.00"
Set-TextCell "E42" "  +0.06%  "
Set-TextCell "D43" "18.13"
Set-TextCell "E43" "  +6.15%  "
Set-TextCell "E44" "  +8.63%  "
Set-TextCell "D45" "1.22"
Set-TextCell "E45" "  +2.91%  "
Set-TextCell "B46" "FTXToken"
Set-TextCell "C46" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D46" "4.46"
Set-TextCell "E46" "  +2.36%  "
Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "98.17"
Set-TextCell "E47" "  +1.26%  "
Set-TextCell "D48" "0.0949"
Set-TextCell "E48" "  +0.61%  "
Set-TextCell "D49" "1.444.28"
Set-TextCell "E49" "  +1.18%  "
Set-TextCell "D50" "2.574.59"
Set-TextCell "E50" "  +5.70%  "
Set-TextCell "E51" "  -2.65%  "
